# Insert a new weekly data row at row 84 (pushing the existing rows 84-124
# down to 85-125) and populate it with the new week's price data for
# Orégano at Mercado Mayorista Lo Valledor de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 84..124 down to 85..125, leaving a blank row 84 to fill in.
$ws.Rows.Item(84).Insert()

$ws.Range("A84").Value = 6
$ws.Range("B84").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C84").Value = "Metropolitana"
$ws.Range("D84").Value = 44523
$ws.Range("E84").Value = 13
$ws.Range("F84").Value = 100112029
$ws.Range("G84").Value = "Orégano"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 35
$ws.Range("K84").Value = 8500
$ws.Range("L84").Value = 9000
$ws.Range("M84").Value = 8729
$ws.Range("N84").Value = "`$/docena de atados"
$ws.Range("O84").Value = "Región Metropolitana"
$ws.Range("P84").Value = 2910
$ws.Range("Q84").Value = 3
$ws.Range("R84").Value = "Hortaliza"
